# Apply updated supply side/IRA assumptions to the "current" run rows
# (Federal Corporate Taxes Contribution, row 6; Fiscal Impact, row 16)
# and recompute the corresponding "difference" rows (34 and 44), which
# are current - previous for each quarterly column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns H..W hold quarterly values from 2022 Q3 through 2026 Q2.
$cols = @("H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W")

# New "current" values for row 6 (Federal Corporate Taxes Contribution)
$row6 = @{
    "H" = 0.0969
    "I" = -0.045
    "J" = 0.2751
    "K" = 0.2384
    "L" = 0.5026
    "M" = 0.0985
    "N" = 0.087
    "O" = 0.0027
    "P" = -0.0924
    "Q" = -0.0876
    "R" = -0.1505
    "S" = -0.0734
    "T" = -0.0804
    "U" = -0.0891
    "V" = -0.0788
    "W" = -1.0948
}

# New "current" values for row 16 (Fiscal Impact)
$row16 = @{
    "H" = -2.2972
    "I" = -0.5605
    "J" = 0.0439
    "K" = -0.2074
    "L" = 0.7225
    "M" = 0.0501
    "N" = -0.3236
    "O" = -0.6367
    "P" = -0.309
    "Q" = -0.0824
    "R" = -0.6647
    "S" = -0.7811
    "T" = -0.5699
    "U" = -0.5089
    "V" = -0.124
    "W" = -72.4415
}

# New "difference" values for row 34 (Federal Corporate Taxes Contribution, current - previous)
$row34 = @{
    "H" = -0.0098
    "I" = 0.0002
    "J" = -0.0199
    "K" = -0.019
    "L" = 0.1938
    "M" = -0.1811
    "N" = -0.0691
    "O" = 0.1229
    "P" = 0.1457
    "Q" = 0.1359
    "R" = 0.1334
    "S" = 0.146
    "T" = 0.0903
    "U" = -0.0537
    "V" = 0.0354
    "W" = -0.6676
}

# New "difference" values for row 44 (Fiscal Impact, current - previous)
$row44 = @{
    "H" = -0.0098
    "I" = 0.0002
    "J" = -0.0199
    "K" = -0.019
    "L" = 0.1938
    "M" = -0.1811
    "N" = -0.0691
    "O" = 0.0371
    "P" = 0.054
    "Q" = 0.0367
    "R" = 0.0215
    "S" = 0.0384
    "T" = -0.0057
    "U" = -0.1384
    "V" = -0.0331
    "W" = -1.7944
}

foreach ($col in $cols) {
    $ws.Range("$col" + "6").Value = $row6[$col]
    $ws.Range("$col" + "16").Value = $row16[$col]
    $ws.Range("$col" + "34").Value = $row34[$col]
    $ws.Range("$col" + "44").Value = $row44[$col]
}
